# "Generate Report for Handback"
#
# The localization-status report previously only recorded the *handoff*
# side of the pipeline (source file + handoff target). This change fills
# in the *handback* side once a handback has actually happened:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - "Latest Target File" (F) / "Latest Handback File" (G) columns are
#     populated with hyperlinked file names (re-using the same files shown
#     in columns A/D, since the handback target is the same file).
#   - "Latest Handback DateTime" (H) is stamped with the real handback
#     timestamp instead of the zero-date placeholder.
#
# This happened for the `zh-cn` and `de-de` locales (and is reflected back
# on the `Overview` sheet, which shares the same "Status" cells).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$aMd = "a.md"
$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/47d5e221aa141db6dcaf0ece3d6903a370720360/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: only the Status text changes (B2/C2/B3/C3).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f899b34c113ed8b70de3a0a8228e131d38d9f83/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Status column
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Populate the new Latest Target File (F) / Latest Handback File (G) cells
$wsZh.Range("F2").Value = $aMd
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("F3").Value = $aMd
$wsZh.Range("G3").Value = $zhXlf

# Match the Hyperlink look used elsewhere on the row (underline + the
# workbook's HyperLink theme color).
$wsZh.Range("F2:G3").Font.Underline = 2
$wsZh.Range("F2:G3").Font.Color = 15570276

# Real handback timestamp replaces the 0001-01-01 placeholder.
$wsZh.Range("H2").Value = "2016-03-25 12:40:50"
$wsZh.Range("H3").Value = "2016-03-25 12:40:50"

# Rebuild the hyperlinks in row order (A2, D2, F2, G2, A3, D3, F3, G3) so
# the relationship ids come out sequential, same as Excel would emit them.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $aMdUrl, "", "", $aMd)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $aMdUrl, "", "", $aMd)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47d5e221aa141db6dcaf0ece3d6903a370720360/e2e/b.md", "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $aMdUrl, "", "", $aMd)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", $zhXlf)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1da9261684e827ac8eb9be2aa674e97490f84ad0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Status column
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Populate the new Latest Target File (F) / Latest Handback File (G) cells
$wsDe.Range("F2").Value = $aMd
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("F3").Value = $aMd
$wsDe.Range("G3").Value = $deXlf

$wsDe.Range("F2:G3").Font.Underline = 2
$wsDe.Range("F2:G3").Font.Color = 15570276

# Real handback timestamp (different from zh-cn's, since de-de handed back
# a minute later) replaces the 0001-01-01 placeholder.
$wsDe.Range("H2").Value = "2016-03-25 12:41:01"
$wsDe.Range("H3").Value = "2016-03-25 12:41:01"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $aMdUrl, "", "", $aMd)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $aMdUrl, "", "", $aMd)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47d5e221aa141db6dcaf0ece3d6903a370720360/e2e/b.md", "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $aMdUrl, "", "", $aMd)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", $deXlf)
